$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.577.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.74"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.509.93"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.83%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.105.35"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.29%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.509.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.570.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.65"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.579"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.650.32"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.04%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.32%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.513.91"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.91%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "24.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.28"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0812"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.814"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.03"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.467.25"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.91"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.04%  "
